$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:09"

# Refresh province/region case counts (rows 4-60), newly sorted by total cases descending;
# also de-duplicates the stray "Cataluna*" / "Cataluna" entry into a single "Cataluna" row.
$ws.Range("A4").Value = "Madrid"
$ws.Range("B4").Value = 6777
$ws.Range("C4").Value = 941
$ws.Range("D4").Value = 5338
$ws.Range("E4").Value = 498
$ws.Range("A5").Value = "Cataluña"
$ws.Range("B5").Value = 3270
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 3185
$ws.Range("E5").Value = 82
$ws.Range("A6").Value = "Araba/Alava"
$ws.Range("B6").Value = 621
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 583
$ws.Range("E6").Value = 35
$ws.Range("A7").Value = "Valencia/Valencia"
$ws.Range("B7").Value = 522
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 503
$ws.Range("E7").Value = 12
$ws.Range("A8").Value = "Navarra"
$ws.Range("B8").Value = 482
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 476
$ws.Range("E8").Value = 4
$ws.Range("A9").Value = "La Rioja"
$ws.Range("B9").Value = 468
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 461
$ws.Range("E9").Value = 5
$ws.Range("A10").Value = "Bizkaia/Vizcaya"
$ws.Range("B10").Value = 393
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = 380
$ws.Range("E10").Value = 13
$ws.Range("A11").Value = "Malaga"
$ws.Range("B11").Value = 361
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 348
$ws.Range("E11").Value = 13
$ws.Range("A12").Value = "Alacant/Alicante"
$ws.Range("B12").Value = 338
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 316
$ws.Range("E12").Value = 11
$ws.Range("A13").Value = "Toledo"
$ws.Range("B13").Value = 293
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 265
$ws.Range("E13").Value = 14
$ws.Range("A14").Value = "Asturias"
$ws.Range("B14").Value = 292
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 287
$ws.Range("E14").Value = 1
$ws.Range("A15").Value = "Albacete"
$ws.Range("B15").Value = 259
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 231
$ws.Range("E15").Value = 20
$ws.Range("A16").Value = "Zaragoza"
$ws.Range("B16").Value = 224
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 210
$ws.Range("E16").Value = 14
$ws.Range("A17").Value = "A Coruña"
$ws.Range("B17").Value = 222
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 219
$ws.Range("E17").Value = 3
$ws.Range("A18").Value = "Ciudad Real"
$ws.Range("B18").Value = 216
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 189
$ws.Range("E18").Value = 21
$ws.Range("A19").Value = "Guadalajara"
$ws.Range("B19").Value = 205
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 3
$ws.Range("A20").Value = "Burgos"
$ws.Range("B20").Value = 187
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 165
$ws.Range("E20").Value = 8
$ws.Range("A21").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B21").Value = 176
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 171
$ws.Range("E21").Value = 5
$ws.Range("A22").Value = "Granada"
$ws.Range("B22").Value = 176
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 169
$ws.Range("E22").Value = 7
$ws.Range("A23").Value = "Aragon"
$ws.Range("B23").Value = 174
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 163
$ws.Range("E23").Value = 11
$ws.Range("A24").Value = "Illes Balears*"
$ws.Range("B24").Value = 169
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 161
$ws.Range("E24").Value = 2
$ws.Range("A25").Value = "Illes Balears"
$ws.Range("B25").Value = 169
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 161
$ws.Range("E25").Value = 2
$ws.Range("A26").Value = "Murcia"
$ws.Range("B26").Value = 168
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 167
$ws.Range("E26").Value = 0
$ws.Range("A27").Value = "Caceres"
$ws.Range("B27").Value = 164
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 153
$ws.Range("E27").Value = 9
$ws.Range("A28").Value = "Salamanca"
$ws.Range("B28").Value = 149
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 129
$ws.Range("E28").Value = 12
$ws.Range("A29").Value = "Pontevedra"
$ws.Range("B29").Value = 145
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 143
$ws.Range("E29").Value = 2
$ws.Range("A30").Value = "Tenerife"
$ws.Range("B30").Value = 143
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 135
$ws.Range("E30").Value = 2
$ws.Range("A31").Value = "Leon"
$ws.Range("B31").Value = 134
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 126
$ws.Range("E31").Value = 5
$ws.Range("A32").Value = "Sevilla"
$ws.Range("B32").Value = 133
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 131
$ws.Range("E32").Value = 1
$ws.Range("A33").Value = "Segovia"
$ws.Range("B33").Value = 121
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 111
$ws.Range("E33").Value = 7
$ws.Range("A34").Value = "Valladolid"
$ws.Range("B34").Value = 115
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 111
$ws.Range("E34").Value = 3
$ws.Range("A35").Value = "Cordoba"
$ws.Range("B35").Value = 101
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 101
$ws.Range("E35").Value = 0
$ws.Range("A36").Value = "Jaen"
$ws.Range("B36").Value = 87
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 85
$ws.Range("E36").Value = 2
$ws.Range("A37").Value = "Cadiz"
$ws.Range("B37").Value = 84
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 84
$ws.Range("E37").Value = 0
$ws.Range("A38").Value = "Cantabria"
$ws.Range("B38").Value = 83
$ws.Range("C38").Value = 10
$ws.Range("D38").Value = 72
$ws.Range("E38").Value = 1
$ws.Range("A39").Value = "Badajoz"
$ws.Range("B39").Value = 77
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 73
$ws.Range("E39").Value = 0
$ws.Range("A40").Value = "Cuenca"
$ws.Range("B40").Value = 72
$ws.Range("C40").Value = 4
$ws.Range("D40").Value = 64
$ws.Range("E40").Value = 4
$ws.Range("A41").Value = "Castello/Castellon"
$ws.Range("B41").Value = 59
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 57
$ws.Range("E41").Value = 1
$ws.Range("A42").Value = "Avila"
$ws.Range("B42").Value = 59
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 55
$ws.Range("E42").Value = 2
$ws.Range("A43").Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Range("B43").Value = 58
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 58
$ws.Range("E43").Value = 3
$ws.Range("A44").Value = "Soria"
$ws.Range("B44").Value = 58
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 49
$ws.Range("E44").Value = 5
$ws.Range("A45").Value = "Gran Canaria"
$ws.Range("B45").Value = 55
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 54
$ws.Range("E45").Value = 1
$ws.Range("A46").Value = "Ourense"
$ws.Range("B46").Value = 46
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = 46
$ws.Range("E46").Value = 0
$ws.Range("A47").Value = "Almeria"
$ws.Range("B47").Value = 37
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 37
$ws.Range("E47").Value = 0
$ws.Range("A48").Value = "Lugo"
$ws.Range("B48").Value = 36
$ws.Range("C48").Value = 4
$ws.Range("D48").Value = 36
$ws.Range("E48").Value = 0
$ws.Range("A49").Value = "Zamora"
$ws.Range("B49").Value = 31
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 29
$ws.Range("E49").Value = 1
$ws.Range("A50").Value = "Teruel"
$ws.Range("B50").Value = 27
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 26
$ws.Range("E50").Value = 1
$ws.Range("A51").Value = "Huesca"
$ws.Range("B51").Value = 24
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 24
$ws.Range("E51").Value = 0
$ws.Range("A52").Value = "Huelva"
$ws.Range("B52").Value = 23
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 23
$ws.Range("E52").Value = 0
$ws.Range("A53").Value = "Melilla"
$ws.Range("B53").Value = 23
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 23
$ws.Range("E53").Value = 0
$ws.Range("A54").Value = "Palencia"
$ws.Range("B54").Value = 14
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 13
$ws.Range("E54").Value = 0
$ws.Range("A55").Value = "Fuerteventura"
$ws.Range("B55").Value = 11
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 11
$ws.Range("E55").Value = 0
$ws.Range("A56").Value = "Arroyo de la Luz"
$ws.Range("B56").Value = 7
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 7
$ws.Range("E56").Value = 0
$ws.Range("A57").Value = "La Palma"
$ws.Range("B57").Value = 5
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 5
$ws.Range("E57").Value = 0
$ws.Range("A58").Value = "Ceuta"
$ws.Range("B58").Value = 5
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 5
$ws.Range("E58").Value = 0
$ws.Range("A59").Value = "Lanzarote"
$ws.Range("B59").Value = 3
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 3
$ws.Range("E59").Value = 0
$ws.Range("A60").Value = "La Gomera"
$ws.Range("B60").Value = 3
$ws.Range("C60").Value = 2
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 0

# The dataset now has one fewer region row; drop the stale trailing row (old row 61)
$ws.Rows.Item(61).Delete()
